# Rewrites the "some_rand(...)" paragraph (collapsing the spell-checked
# runs into a single run) and adds two new documentation paragraphs plus
# a blank spacer paragraph, describing self.subscription_pool and
# self.subscription_queue.

$d = $word.ActiveDocument

function Insert-Ooxml($rng, [string]$bodyXml) {
    # Wraps a <w:p>...</w:p> fragment (or several) in a minimal OPC
    # package so Range.InsertXML can splice it in at $rng.
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg) | Out-Null
}

# Smart/curly quotes used in the subscription_pool comment.
$ldq = [char]0x201C
$rdq = [char]0x201D

# --- Locate the existing "some_rand(...)" paragraph (3rd paragraph). ---
$target = $d.Paragraphs.Item(3)
$rng = $target.Range
# Exclude the trailing paragraph mark so we only clear the run content.
$rng.End = $rng.End - 1
$rng.Delete()

# 1) Replace the cleared paragraph with a single run containing the tab
#    and the merged (no longer spell-split) sentence.
$rng1 = $d.Paragraphs.Item(3).Range
$rng1.Collapse(1)
Insert-Ooxml $rng1 '<w:p><w:r><w:tab/><w:t>some_rand(x, a, b): generate x different ints in [a,b)</w:t></w:r></w:p>'

# 2) Blank spacer paragraph.
$rng2 = $d.Paragraphs.Item(4).Range
$rng2.Collapse(1)
Insert-Ooxml $rng2 '<w:p/>'

# 3) New paragraph documenting self.subscription_pool.
$rng3 = $d.Paragraphs.Item(5).Range
$rng3.Collapse(1)
$body3 = '<w:p>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>self</w:t></w:r>' +
    '<w:r><w:t>.subscription_pool: save all subscribe info</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> [{' + $ldq + 'topic name' + $rdq + ': [ broker names]}]</w:t></w:r>' +
    '</w:p>'
Insert-Ooxml $rng3 $body3

# 4) New paragraph documenting self.subscription_queue.
$rng4 = $d.Paragraphs.Item(6).Range
$rng4.Collapse(1)
$body4 = '<w:p>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:tab/><w:t>self.subscription_queue: subscription info that need to be flood [(topic, name), ]</w:t></w:r>' +
    '</w:p>'
Insert-Ooxml $rng4 $body4
